# add 3 groups of data hyy hzj cxq
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings. Write in this order so the shared-string table picks
# up new unique values in the same sequence as the target workbook
# (cxq, then hyy, then hzj).
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

# New numeric data for the 3 new columns
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 0.98250728862973757

$ws.Range("F2").Value = 0.98498498498498499
$ws.Range("F3").Value = 0.97741935483870968

$ws.Range("G2").Value = 0.97979797979797978
$ws.Range("G3").Value = 0.95221843003412965

# Widen the new last column
$ws.Range("G1").ColumnWidth = 20.5

# Leave the selection where it ends up after entering this data
$ws.Range("G10").Select()
